$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Valor" -> shifts to C)
$ws.Columns("B").Insert()

# Header row
$ws.Range("B1").Value = "Variável"
$ws.Range("C1").Value = "Valor"
$ws.Range("D1").Value = "Colocação"

# Copy header style (bold/centered) from A1 to the new header cells
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Column B: "Variação 2021/2012" for rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "Variação 2021/2012"
}

# Column D: Colocação (ranking) for rows 2-8
$rankings = @("1º", "2º", "3º", "4º", "5º", "6º", "21º")
for ($i = 0; $i -lt $rankings.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 4).Value = $rankings[$i]
}
